$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in Code ID column (D) for rows 2-65 with sequential values 401-464
for ($row = 2; $row -le 65; $row++) {
    $ws.Cells.Item($row, 4).Value = 400 + ($row - 1)
}

# Update the visible selection to match the post-edit state (C2:C65, active C2)
$ws.Range("C2:C65").Select()
